$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - updated TPM-derived values
$ws.Range("M2").Value = 4.621579
$ws.Range("N2").Value = 13.864737
$ws.Range("O2").Value = 0.1778708528171788
$ws.Range("P2").Value = 0.1778708528171788
$ws.Range("Q2").Value = 0.140514487916
$ws.Range("R2").Value = 1.264630391244
$ws.Range("S2").Value = 0.1778708528171788
$ws.Range("T2").Value = 0.1778708528171788

# Row 3 - updated TPM-derived values
$ws.Range("N3").Value = 46.543441
$ws.Range("O3").Value = 0.5971062807549863
$ws.Range("P3").Value = 0.5971062807549863
$ws.Range("R3").Value = 4.245320340492
$ws.Range("S3").Value = 0.5971062807549863
$ws.Range("T3").Value = 0.5971062807549863

# Row 4 - updated TPM-derived values
$ws.Range("O4").Value = 0.2250228664278349
$ws.Range("P4").Value = 0.2250228664278349
$ws.Range("S4").Value = 0.2250228664278349
$ws.Range("T4").Value = 0.2250228664278349
